# Revised marking sheet for Jungah's thesis:
# - Criterion 4 ("Discussion/analysis...", row 6 of the Grading sheet) is
#   revised from a 5 down to a 4, which drops the overall weighted grade
#   (formula in E11) from "A" down to "B".
# - The grader's cursor/selection on the Grading sheet ends up on E9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grading")
$ws.Activate()

# Lower the score for the 4th grading criterion from 5 to 4.
$ws.Range("E6").Value = 4

# Scroll the window down a bit and leave the selection on E9, matching
# where the grader was last working.
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E9").Select()
